$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 257, pushing the existing rows 257-285
# down to become rows 259-287 (dimension grows from A1:R285 to A1:R287).
$ws.Rows.Item(257).Insert()
$ws.Rows.Item(257).Insert()

# Populate the first new row (257) with its data.
$ws.Cells.Item(257,1).Value = 7
$ws.Cells.Item(257,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(257,3).Value = "Ñuble"
$ws.Cells.Item(257,4).Value = 44769
$ws.Cells.Item(257,5).Value = 16
$ws.Cells.Item(257,6).Value = 100112023
$ws.Cells.Item(257,7).Value = "Brócoli"
$ws.Cells.Item(257,8).Value = "Sin especificar"
$ws.Cells.Item(257,9).Value = "Primera"
$ws.Cells.Item(257,10).Value = 200
$ws.Cells.Item(257,11).Value = 900
$ws.Cells.Item(257,12).Value = 1000
$ws.Cells.Item(257,13).Value = 950
$ws.Cells.Item(257,14).Value = "`$/unidad"
$ws.Cells.Item(257,15).Value = "Provincia de Diguillín"
$ws.Cells.Item(257,16).Value = 950
$ws.Cells.Item(257,17).Value = 1
$ws.Cells.Item(257,18).Value = "Hortaliza"

# Populate the second new row (258) with its data.
$ws.Cells.Item(258,1).Value = 7
$ws.Cells.Item(258,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(258,3).Value = "Ñuble"
$ws.Cells.Item(258,4).Value = 44769
$ws.Cells.Item(258,5).Value = 16
$ws.Cells.Item(258,6).Value = 100112023
$ws.Cells.Item(258,7).Value = "Brócoli"
$ws.Cells.Item(258,8).Value = "Sin especificar"
$ws.Cells.Item(258,9).Value = "Segunda"
$ws.Cells.Item(258,10).Value = 100
$ws.Cells.Item(258,11).Value = 800
$ws.Cells.Item(258,12).Value = 800
$ws.Cells.Item(258,13).Value = 800
$ws.Cells.Item(258,14).Value = "`$/unidad"
$ws.Cells.Item(258,15).Value = "Provincia de Diguillín"
$ws.Cells.Item(258,16).Value = 800
$ws.Cells.Item(258,17).Value = 1
$ws.Cells.Item(258,18).Value = "Hortaliza"
